$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update its "through" date label (2022-10-31 -> 2022-11-01)
$ws.Name = "Through 2022-11-01"
$ws.Range("I1").Value = "2022 (through 11-01)"

# Add new data for November (row 12), column I (new year column)
$ws.Range("I12").Value = 2

# Update the yearly total for column I (row 14)
$ws.Range("I14").Value = 1403
